$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row, bold style already applied to row) ---
$ws.Range("A1").Value = "Status"
$ws.Range("B1").Value = "Issue"
$ws.Range("C1").Value = "Comment"

# --- Row 2 ---
$ws.Range("A2").Value = "Acceptable"
$ws.Range("B2").Value = "OAL comments are not preserved"
$ws.Range("C2").Value = "These do not appear to be in the OAL metamodel"

# --- Row 3 ---
$ws.Range("B3").Value = "declare/begin/end not generated for any statement type"
$ws.Range("C3").Value = "Cannot find where to do this"

# --- Row 4 ---
$ws.Range("B4").Value = "Event generation parameters are not generated in the correct order"
$ws.Range("C4").Value = "These should be in the same order as the model declaration.  Tried suggestions from Cort but no success."

# --- Row 5 ---
$ws.Range("B5").Value = "Operation call parameters are not generated in the correct order"
$ws.Range("C5").Value = "These should be in the same order as the model declaration.  Tried suggestions from Cort but no success."

# --- Row 6 ---
$ws.Range("A6").Value = "Acceptable"
$ws.Range("B6").Value = "Multiple elifs with no else not supported"
$ws.Range("C6").Value = "This is due to the way that the current translation structure utilises the C/C++ block syntax"

# --- Row 7 ---
$ws.Range("B7").Value = 'Select where conditions have "SELECTED." and "=="'
$ws.Range("C7").Value = 'Have got rid of SLECTED but "==" and some extraneous "." occurrences remain'

# --- Row 8 ---
$ws.Range("B8").Value = "Enumeration literals not rendered correctly"
$ws.Range("C8").Value = "Should be just the literal name (not the mangled C/C++ name)"

# --- Row 9 ---
$ws.Range("B9").Value = "Variable declarations are not correct"
$ws.Range("C9").Value = "Remnants of C/C++ declarations"

# --- Row 10 ---
$ws.Range("B10").Value = "Comment has been inserted for the create event instance statement but not for the bridge call to the create timer"
$ws.Range("C10").Value = "Not a problem as a parse error and comment has been inserted and use will have to manually sort this out anyway"

# --- Row 11 ---
$ws.Range("A11").Value = "Fixed"
$ws.Range("B11").Value = 'State action and service/function bodies reference "PARAM."'
$ws.Range("C11").Value = "Not required in MASL"

# --- Row 12 ---
$ws.Range("B12").Value = "Function calls as part of an expression not translated"
$ws.Range("C12").Value = "Cannot find where to do this"

# --- Row 13 ---
$ws.Range("A13").Value = "Acceptable"
$ws.Range("B13").Value = '"Control" statement not implemented (Parse error and comment is generated)'
$ws.Range("C13").Value = "This allows architecture dependent calls.  Not relevant to MASL."

# --- Row 14 ---
$ws.Range("A14").Value = "Acceptable"
$ws.Range("B14").Value = '"Continue" statement not implemented (Parse error and comment is generated)'
$ws.Range("C14").Value = 'There appears to be no "continue" in MASL.  Implementing this by elaborated MASL would be very complex'

# --- Row 15 ---
$ws.Range("B15").Value = "Relationship navigation not yet tested for associative classes and super/sub hierarchies"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.1666666666666667
$ws.Columns.Item(2).ColumnWidth = 94.60807291666667

# --- Selection ---
$ws.Range("B3").Select() | Out-Null

# --- Page setup: landscape, 65% scale, fit to 1 page tall ---
$ws.PageSetup.Zoom = 65
$ws.PageSetup.FitToPagesTall = $false
$ws.PageSetup.Orientation = 2
